# Snell's Law workbook — "Further bug fix and cosmetic enhancement"
#
# 1) Move the incident-angle scroll bar to 79 degrees (drives the whole
#    calculation cascade) and correct the outgoing medium's refractive
#    index to 1.33 (was 1).
# 2) Fix the outgoing-angle formula (Calculations!B7) so that the
#    90/270-degree edge cases return the angle itself instead of relying
#    on the (buggy) generic branches.
# 3) Cosmetic: move the active selection on "Inputs & Outputs" to J15.

$wb = $excel.ActiveWorkbook

$wsIO   = $wb.Worksheets.Item("Inputs & Outputs")
$wsCalc = $wb.Worksheets.Item("Calculations")

# --- Data changes -----------------------------------------------------

# Incident angle (linked to the "Scroll Bar 1" form control).
$wsIO.Range("B2").Value = 79

# Keep the scroll bar control itself in sync with its linked cell.
$scrollBar = $wsIO.Shapes.Item("Scroll Bar 1")
$cf = $scrollBar.ControlFormat
$cf.Min = 0
$cf.Max = 360
$cf.SmallChange = 1
$cf.LargeChange = 10
$cf.Value = 79

# Refractive index of the outgoing medium.
$wsIO.Range("B5").Value = 1.33

# --- Bug fix: Calculations!B7 -----------------------------------------
# Handle the D2 = 90 / D2 = 270 edge cases explicitly instead of falling
# through to the generic branches.
$wsCalc.Range("B7").Formula = "=IF(OR(D2=90,D2=270),D2,IF(D2>270,B6+360,IF(D2>180,180-B6,IF(D2>90,180-B6,B6))))"

# --- Cosmetic: selection -------------------------------------------------
$wsIO.Activate() | Out-Null
$wsIO.Range("J15").Select() | Out-Null
